# queries de tudo ou nada - categoria-menção
# Corrige a troca de rótulos de UF nas linhas finais (empate / categoria "menção")
# das planilhas de total arrecadado e máximo arrecadado.

$wb = $excel.ActiveWorkbook

# uf-tot-arrecad: A25/A26 estavam TO/AC -> devem ficar AC/TO
$wsTot = $wb.Worksheets.Item("uf-tot-arrecad")
$wsTot.Range("A25").Value = "AC"
$wsTot.Range("A26").Value = "TO"

# uf-max-arrecad: A25/A26 estavam AC/TO -> devem ficar TO/AC
$wsMax = $wb.Worksheets.Item("uf-max-arrecad")
$wsMax.Range("A25").Value = "TO"
$wsMax.Range("A26").Value = "AC"
